$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.551609188260556
$ws.Cells.Item(2, 3).Value = 0.2002258939387218
$ws.Cells.Item(2, 5).Value = 0.1693266460553522
$ws.Cells.Item(2, 6).Value = 2.342751699717581
$ws.Cells.Item(2, 7).Value = 1.317092232050996
$ws.Cells.Item(2, 8).Value = 1.220246611047443
$ws.Cells.Item(2, 10).Value = 0.06949057467695852
$ws.Cells.Item(2, 12).Value = 0.3766328660746296
$ws.Cells.Item(2, 13).Value = 0.3836746348334543
$ws.Cells.Item(2, 14).Value = 1.794068942494967

$ws.Cells.Item(3, 2).Value = 1.461575394362455
$ws.Cells.Item(3, 3).Value = 0.1852891020234892
$ws.Cells.Item(3, 5).Value = 0.1696996610878445
$ws.Cells.Item(3, 6).Value = 2.338045565844794
$ws.Cells.Item(3, 7).Value = 1.307905292019186
$ws.Cells.Item(3, 8).Value = 1.221959196064518
$ws.Cells.Item(3, 10).Value = 0.06979986181996445
$ws.Cells.Item(3, 12).Value = 0.3735828652350861
$ws.Cells.Item(3, 13).Value = 0.369668767803546
$ws.Cells.Item(3, 14).Value = 1.815525494702882

$ws.Cells.Item(4, 2).Value = 1.406979468205748
$ws.Cells.Item(4, 3).Value = 0.176034563525505
$ws.Cells.Item(4, 5).Value = 0.1699477070854929
$ws.Cells.Item(4, 6).Value = 2.336379663702147
$ws.Cells.Item(4, 7).Value = 1.303165173775909
$ws.Cells.Item(4, 8).Value = 1.223647629460103
$ws.Cells.Item(4, 10).Value = 0.07000202509476239
$ws.Cells.Item(4, 12).Value = 0.3718354269731492
$ws.Cells.Item(4, 13).Value = 0.3612399859500286
$ws.Cells.Item(4, 14).Value = 1.829352119404404

$ws.Cells.Item(5, 2).Value = 1.384903980386525
$ws.Cells.Item(5, 3).Value = 0.1722422475864249
$ws.Cells.Item(5, 5).Value = 0.1700535805241925
$ws.Cells.Item(5, 6).Value = 2.336008392945558
$ws.Cells.Item(5, 7).Value = 1.301459438936419
$ws.Cells.Item(5, 8).Value = 1.224495697345134
$ws.Cells.Item(5, 10).Value = 0.07008750032596289
$ws.Cells.Item(5, 12).Value = 0.3711549555423019
$ws.Cells.Item(5, 13).Value = 0.3578483045246017
$ws.Cells.Item(5, 14).Value = 1.835150484476123

$ws.Cells.Item(6, 2).Value = 1.381248813856189
$ws.Cells.Item(6, 3).Value = 0.1716112649808821
$ws.Cells.Item(6, 5).Value = 0.1700714505386922
$ws.Cells.Item(6, 6).Value = 2.335965318851123
$ws.Cells.Item(6, 7).Value = 1.301189829800222
$ws.Cells.Item(6, 8).Value = 1.224646179119304
$ws.Cells.Item(6, 10).Value = 0.07010188045005705
$ws.Cells.Item(6, 12).Value = 0.3710438773991029
$ws.Cells.Item(6, 13).Value = 0.3572877268603776
$ws.Cells.Item(6, 14).Value = 1.836123195034473

$ws.Cells.Item(7, 2).Value = 1.40668105014953
$ws.Cells.Item(7, 3).Value = 0.1759835042143436
$ws.Cells.Item(7, 5).Value = 0.1699491155100676
$ws.Cells.Item(7, 6).Value = 2.33637341130968
$ws.Cells.Item(7, 7).Value = 1.303141255751527
$ws.Cells.Item(7, 8).Value = 1.223658419080309
$ws.Cells.Item(7, 10).Value = 0.07000316531221529
$ws.Cells.Item(7, 12).Value = 0.3718261216895584
$ws.Cells.Item(7, 13).Value = 0.3611940697787759
$ws.Cells.Item(7, 14).Value = 1.829429654685593

$ws.Cells.Item(8, 2).Value = 1.52042359153802
$ws.Cells.Item(8, 3).Value = 0.1950929263774981
$ws.Cells.Item(8, 5).Value = 0.169451322391641
$ws.Cells.Item(8, 6).Value = 2.340874976139531
$ws.Cells.Item(8, 7).Value = 1.31373727438951
$ws.Cells.Item(8, 8).Value = 1.220704836759893
$ws.Cells.Item(8, 10).Value = 0.06959467789649665
$ws.Cells.Item(8, 12).Value = 0.375555296749404
$ws.Cells.Item(8, 13).Value = 0.378810028992369
$ws.Cells.Item(8, 14).Value = 1.801331735132801

$ws.Cells.Item(9, 2).Value = 1.748900358376488
$ws.Cells.Item(9, 3).Value = 0.2319105352412976
$ws.Cells.Item(9, 5).Value = 0.1686255054868524
$ws.Cells.Item(9, 6).Value = 2.359420467312688
$ws.Cells.Item(9, 7).Value = 1.341692241379405
$ws.Cells.Item(9, 8).Value = 1.21997390665021
$ws.Cells.Item(9, 10).Value = 0.06889050337120572
$ws.Cells.Item(9, 12).Value = 0.3838574999493147
$ws.Cells.Item(9, 13).Value = 0.4147064508793719
$ws.Cells.Item(9, 14).Value = 1.751410528732236

$ws.Cells.Item(10, 2).Value = 1.92007930045213
$ws.Cells.Item(10, 3).Value = 0.2585697925951536
$ws.Cells.Item(10, 5).Value = 0.1681097603765931
$ws.Cells.Item(10, 6).Value = 2.378987394098957
$ws.Cells.Item(10, 7).Value = 1.366650950188983
$ws.Cells.Item(10, 8).Value = 1.222534673083743
$ws.Cells.Item(10, 10).Value = 0.06843165060741541
$ws.Cells.Item(10, 12).Value = 0.3905549379012996
$ws.Cells.Item(10, 13).Value = 0.4419010658977101
$ws.Cells.Item(10, 14).Value = 1.717894360184625

$ws.Cells.Item(11, 2).Value = 1.998676387863497
$ws.Cells.Item(11, 3).Value = 0.2706150876745426
$ws.Cells.Item(11, 5).Value = 0.1678947499932022
$ws.Cells.Item(11, 6).Value = 2.389183178575976
$ws.Cells.Item(11, 7).Value = 1.378975224721017
$ws.Cells.Item(11, 8).Value = 1.22437514754219
$ws.Cells.Item(11, 10).Value = 0.06823549727453226
$ws.Cells.Item(11, 12).Value = 0.3937305363765802
$ws.Cells.Item(11, 13).Value = 0.4544506977105769
$ws.Cells.Item(11, 14).Value = 1.703334221497235

$ws.Cells.Item(12, 2).Value = 2.028543426254657
$ws.Cells.Item(12, 3).Value = 0.2751646352312207
$ws.Cells.Item(12, 5).Value = 0.1678161392287807
$ws.Cells.Item(12, 6).Value = 2.393230475569482
$ws.Cells.Item(12, 7).Value = 1.383782393011018
$ws.Cells.Item(12, 8).Value = 1.225169437938746
$ws.Cells.Item(12, 10).Value = 0.06816301974344086
$ws.Cells.Item(12, 12).Value = 0.3949514817624618
$ws.Cells.Item(12, 13).Value = 0.4592285156882028
$ws.Cells.Item(12, 14).Value = 1.697919525712236

$ws.Cells.Item(13, 2).Value = 2.022106406439491
$ws.Cells.Item(13, 3).Value = 0.2741853300220214
$ws.Cells.Item(13, 5).Value = 0.1678329446921039
$ws.Cells.Item(13, 6).Value = 2.392350525459449
$ws.Cells.Item(13, 7).Value = 1.382740835439819
$ws.Cells.Item(13, 8).Value = 1.224994040783997
$ws.Cells.Item(13, 10).Value = 0.06817854907511212
$ws.Cells.Item(13, 12).Value = 0.3946877122204029
$ws.Cells.Item(13, 13).Value = 0.4581983933136016
$ws.Cells.Item(13, 14).Value = 1.699081275076918

$ws.Cells.Item(14, 2).Value = 2.001131485082681
$ws.Cells.Item(14, 3).Value = 0.2709896166777241
$ws.Cells.Item(14, 5).Value = 0.1678882263990884
$ws.Cells.Item(14, 6).Value = 2.38951241575198
$ws.Cells.Item(14, 7).Value = 1.379367898927853
$ws.Cells.Item(14, 8).Value = 1.224438542262163
$ws.Cells.Item(14, 10).Value = 0.06822949844323567
$ws.Cells.Item(14, 12).Value = 0.3938306159086551
$ws.Cells.Item(14, 13).Value = 0.4548432605522024
$ws.Cells.Item(14, 14).Value = 1.70288676775128

$ws.Cells.Item(15, 2).Value = 1.988297272714306
$ws.Cells.Item(15, 3).Value = 0.2690306225234735
$ws.Cells.Item(15, 5).Value = 0.167922453558135
$ws.Cells.Item(15, 6).Value = 2.38779826908447
$ws.Cells.Item(15, 7).Value = 1.377320160826827
$ws.Cells.Item(15, 8).Value = 1.224110966152921
$ws.Cells.Item(15, 10).Value = 0.06826094078439837
$ws.Cells.Item(15, 12).Value = 0.3933080143024057
$ws.Cells.Item(15, 13).Value = 0.452791467205806
$ws.Cells.Item(15, 14).Value = 1.705230628980881

$ws.Cells.Item(16, 2).Value = 1.914957380287717
$ws.Cells.Item(16, 3).Value = 0.2577809571438081
$ws.Cells.Item(16, 5).Value = 0.1681242054375796
$ws.Cells.Item(16, 6).Value = 2.37834714508287
$ws.Cells.Item(16, 7).Value = 1.365865115994524
$ws.Cells.Item(16, 8).Value = 1.222428004668103
$ws.Cells.Item(16, 10).Value = 0.0684447221927833
$ws.Cells.Item(16, 12).Value = 0.3903499885551298
$ws.Cells.Item(16, 13).Value = 0.441084500535041
$ws.Cells.Item(16, 14).Value = 1.718859730430873

$ws.Cells.Item(17, 2).Value = 1.870151553838639
$ws.Cells.Item(17, 3).Value = 0.2508586622064968
$ws.Cells.Item(17, 5).Value = 0.1682529876490797
$ws.Cells.Item(17, 6).Value = 2.372880928987257
$ws.Cells.Item(17, 7).Value = 1.35908685489764
$ws.Cells.Item(17, 8).Value = 1.221568730740529
$ws.Cells.Item(17, 10).Value = 0.06856068293068063
$ws.Cells.Item(17, 12).Value = 0.3885682640310222
$ws.Cells.Item(17, 13).Value = 0.4339483109322941
$ws.Cells.Item(17, 14).Value = 1.727396714112963

$ws.Cells.Item(18, 2).Value = 1.844448879004517
$ws.Cells.Item(18, 3).Value = 0.2468694062994246
$ws.Cells.Item(18, 5).Value = 0.1683289056511539
$ws.Cells.Item(18, 6).Value = 2.36985875634123
$ws.Cells.Item(18, 7).Value = 1.355279478834319
$ws.Cells.Item(18, 8).Value = 1.221138074885772
$ws.Cells.Item(18, 10).Value = 0.06862856511474824
$ws.Cells.Item(18, 12).Value = 0.3875556030842375
$ws.Cells.Item(18, 13).Value = 0.4298605889037148
$ws.Cells.Item(18, 14).Value = 1.732371569483573

$ws.Cells.Item(19, 2).Value = 1.835758167061385
$ws.Cells.Item(19, 3).Value = 0.2455173822918937
$ws.Cells.Item(19, 5).Value = 0.1683549275244101
$ws.Cells.Item(19, 6).Value = 2.368856422503512
$ws.Cells.Item(19, 7).Value = 1.354006026831144
$ws.Cells.Item(19, 8).Value = 1.221003175479979
$ws.Cells.Item(19, 10).Value = 0.06865175255560452
$ws.Cells.Item(19, 12).Value = 0.3872148222077527
$ws.Cells.Item(19, 13).Value = 0.4284794507627367
$ws.Cells.Item(19, 14).Value = 1.734067060856198

$ws.Cells.Item(20, 2).Value = 1.874914129763852
$ws.Cells.Item(20, 3).Value = 0.2515963523823643
$ws.Cells.Item(20, 5).Value = 0.1682390875956834
$ws.Cells.Item(20, 6).Value = 2.373450204534805
$ws.Cells.Item(20, 7).Value = 1.359798958148957
$ws.Cells.Item(20, 8).Value = 1.22165362063987
$ws.Cells.Item(20, 10).Value = 0.06854821616242823
$ws.Cells.Item(20, 12).Value = 0.3887566760933368
$ws.Cells.Item(20, 13).Value = 0.4347062301936262
$ws.Cells.Item(20, 14).Value = 1.726481249557532

$ws.Cells.Item(21, 2).Value = 2.007289508629867
$ws.Cells.Item(21, 3).Value = 0.2719285933537776
$ws.Cells.Item(21, 5).Value = 0.1678719126661221
$ws.Cells.Item(21, 6).Value = 2.390340977603486
$ws.Cells.Item(21, 7).Value = 1.380354800771329
$ws.Cells.Item(21, 8).Value = 1.224599062209364
$ws.Cells.Item(21, 10).Value = 0.06821448454711465
$ws.Cells.Item(21, 12).Value = 0.3940818670470492
$ws.Cells.Item(21, 13).Value = 0.4558280522597542
$ws.Cells.Item(21, 14).Value = 1.701766316302303

$ws.Cells.Item(22, 2).Value = 2.094410593620864
$ws.Cells.Item(22, 3).Value = 0.2851484651847045
$ws.Cells.Item(22, 5).Value = 0.1676483110486
$ws.Cells.Item(22, 6).Value = 2.402466503999918
$ws.Cells.Item(22, 7).Value = 1.394606938759665
$ws.Cells.Item(22, 8).Value = 1.227091559476406
$ws.Cells.Item(22, 10).Value = 0.06800686863788563
$ws.Cells.Item(22, 12).Value = 0.3976694694207907
$ws.Cells.Item(22, 13).Value = 0.4697812114105773
$ws.Cells.Item(22, 14).Value = 1.686190309785629

$ws.Cells.Item(23, 2).Value = 2.047857134889398
$ws.Cells.Item(23, 3).Value = 0.2780990092350919
$ws.Cells.Item(23, 5).Value = 0.1677661570312683
$ws.Cells.Item(23, 6).Value = 2.395895406988728
$ws.Cells.Item(23, 7).Value = 1.386925256561568
$ws.Cells.Item(23, 8).Value = 1.225709276489681
$ws.Cells.Item(23, 10).Value = 0.06811671915865603
$ws.Cells.Item(23, 12).Value = 0.3957449218182205
$ws.Cells.Item(23, 13).Value = 0.4623205782111413
$ws.Cells.Item(23, 14).Value = 1.694450691691826

$ws.Cells.Item(24, 2).Value = 1.872760792247448
$ws.Cells.Item(24, 3).Value = 0.2512628723694093
$ws.Cells.Item(24, 5).Value = 0.1682453659559009
$ws.Cells.Item(24, 6).Value = 2.373192459939176
$ws.Cells.Item(24, 7).Value = 1.359476737413019
$ws.Cells.Item(24, 8).Value = 1.221615044580176
$ws.Cells.Item(24, 10).Value = 0.06855384860288716
$ws.Cells.Item(24, 12).Value = 0.3886714586224116
$ws.Cells.Item(24, 13).Value = 0.4343635282590981
$ws.Cells.Item(24, 14).Value = 1.726894922911135

$ws.Cells.Item(25, 2).Value = 1.686509441947749
$ws.Cells.Item(25, 3).Value = 0.2220196664868297
$ws.Cells.Item(25, 5).Value = 0.168832884995203
$ws.Cells.Item(25, 6).Value = 2.353361281265265
$ws.Cells.Item(25, 7).Value = 1.333356743509938
$ws.Cells.Item(25, 8).Value = 1.219628549228617
$ws.Cells.Item(25, 10).Value = 0.06907068841438857
$ws.Cells.Item(25, 12).Value = 0.3815061726852065
$ws.Cells.Item(25, 13).Value = 0.4048510807701149
$ws.Cells.Item(25, 14).Value = 1.764360741962646
